$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values after repulling data / recalculating mean.
$updates = @{
    4  = 0
    5  = 0
    13 = 0
    20 = 1
    21 = 2
    27 = -1
    35 = -1
    40 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
